$wb = $excel.ActiveWorkbook

# The "Mahesh" first name entered in Sheet2!A2 was a negative-case test
# value that needs to be removed (its shared string also drops out since
# nothing else references it); the cell keeps its existing style.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A2").ClearContents()

# Sheet2 becomes the active/selected tab (Sheet1 was active before).
$ws2.Select()
